$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Clean up the two grammar-check ("proofErr") run splits so each sentence
#    becomes a single run again (no visible text change, just de-fragmenting
#    runs that Word's grammar checker had split around "mode").
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("used. If mode address", $true, $false, $false, $false, $false, $true, 1, $false, "used. If mode address", 2) | Out-Null
$d.Content.Find.Execute("with the mean of", $true, $false, $false, $false, $false, $true, 1, $false, "with the mean of", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Append the new "4. Feature Scaling:" section at the end of the document.
# ---------------------------------------------------------------------------

# Create the 4 new (still empty / plain) paragraphs up front, while the
# "insertion point" paragraph still carries only the plain sz=24 formatting
# from the preceding paragraph -- this avoids the new paragraphs picking up
# heavy heading formatting later applied to the heading paragraph.
$paras = $d.Paragraphs
$last = $paras.Item($paras.Count)
$last.Range.InsertParagraphAfter()

$paras = $d.Paragraphs
$pSpacer = $paras.Item($paras.Count)
$pSpacer.Range.InsertParagraphAfter()

$paras = $d.Paragraphs
$pHeadingPara = $paras.Item($paras.Count)
$pHeadingPara.Range.InsertParagraphAfter()

$paras = $d.Paragraphs
$pBodyPara = $paras.Item($paras.Count)
$pBodyPara.Range.InsertParagraphAfter()

# Now there are 4 fresh empty paragraphs at the tail of the document:
#   pHeading (index count-2), pBody (index count-1), pTail (index count)
$paras = $d.Paragraphs
$n = $paras.Count
$pHeading = $paras.Item($n - 2)
$pBody = $paras.Item($n - 1)
$pTail = $paras.Item($n)

# --- Heading paragraph: "4. Feature Scaling:" (bold, red, 20pt) -----------
$r = $pHeading.Range
$r.Text = "4. Feature Scaling:"
$r.Font.Bold = $true
$r.Font.BoldBi = $true
$r.Font.Color = 255
$r.Font.Size = 20
$r.Font.SizeBi = 20

# --- Body paragraph: description text (plain, 12pt) ------------------------
$text4 = "To handle highly varying magnitudes or values or units, normalization was performed features which have a relatively small range of values ('category', 'amenities', 'bathrooms','bedrooms','has_photo','pets_allowed','price','square_feet','address','cityname','state','latitude','longitude','source') and standardization on the features with a wide range of values ('title', 'body') to avoid weighing greater values, higher and consider smaller values as the lower values, regardless of the unit of the values"
$pBody.Range.Text = $text4

# --- Trailing empty paragraph with a left tab stop at 1760 twips ----------
$pTail.Range.ParagraphFormat.TabStops.Add(88) | Out-Null

Write-Host "edit complete"
